# "added tabs for RES figures"
# The NPV figures in column B were plain numbers (1967 / 1552 / 1414); they
# are now shown as billion-dollar amount strings, column B is widened to fit
# the new text, and the last active selection moves to D9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "$1967 Billion"
$ws.Range("B3").Value = "$1552 Billion"
$ws.Range("B4").Value = "$1414 Billion"

# Widen column B so the new "$#### Billion" labels fit.
$ws.Columns.Item(2).ColumnWidth = 15.6666666666667

# Leave the selection where the author last left it before saving.
$ws.Range("D9").Select()
